$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking price values so they are not
# auto-converted to numbers by Excel (matches original inlineStr text cells).
$textCells = @("D5","D6","D10","D11","D13","D14","D18","D20","D21","D23","D24","D25","D28","D29","D30","D31","D34","D37","D41","D42","D43","D45","D46","D49")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '65.844.16'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '2.695.18'
$ws.Range('E3').Value = '  +2.32%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '610.09'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').Value = '158.18'
$ws.Range('E6').Value = '  +2.24%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.39%  '
$ws.Range('E9').Value = '  +6.97%  '
$ws.Range('D10').Value = '6.03'
$ws.Range('E10').Value = '  +4.72%  '
$ws.Range('D11').Value = '0.403'
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('D13').Value = '30.31'
$ws.Range('E13').Value = '  +5.17%  '
$ws.Range('D14').Value = '0.0000202'
$ws.Range('E14').Value = '  +9.38%  '
$ws.Range('D15').Value = '3.183.09'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').Value = '65.713.49'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').Value = '2.704.03'
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').Value = '12.63'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('E19').Value = '  +2.70%  '
$ws.Range('D20').Value = '359.71'
$ws.Range('E20').Value = '  +2.88%  '
$ws.Range('D21').Value = '7.55'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '70.80'
$ws.Range('E23').Value = '  +4.52%  '
$ws.Range('D24').Value = '9.87'
$ws.Range('E24').Value = '  +4.46%  '
$ws.Range('D25').Value = '0.0000107'
$ws.Range('E25').Value = '  +15.30%  '
$ws.Range('E26').Value = '  -1.64%  '
$ws.Range('E27').Value = '  +3.58%  '
$ws.Range('D28').Value = '0.173'
$ws.Range('E28').Value = '  +6.14%  '
$ws.Range('D29').Value = '8.34'
$ws.Range('E29').Value = '  +3.22%  '
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').Value = '546.92'
$ws.Range('E30').Value = '  +7.36%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '2.21'
$ws.Range('E31').Value = '  +6.63%  '
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('E33').Value = '  +2.12%  '
$ws.Range('D34').Value = '6.69'
$ws.Range('E34').Value = '  +8.09%  '
$ws.Range('E35').Value = '  -3.70%  '
$ws.Range('E36').Value = '  +2.54%  '
$ws.Range('D37').Value = '20.83'
$ws.Range('E37').Value = '  +3.95%  '
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '170.79'
$ws.Range('E41').Value = '  +4.48%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '42.91'
$ws.Range('E43').Value = '  +1.76%  '
$ws.Range('E44').Value = '  +3.19%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '23.66'
$ws.Range('E45').Value = '  +4.39%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0614'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('E47').Value = '  +4.91%  '
$ws.Range('E48').Value = '  +5.67%  '
$ws.Range('D49').Value = '0.660'
$ws.Range('E49').Value = '  +2.62%  '
$ws.Range('E50').Value = '  +9.69%  '
$ws.Range('E51').Value = '  +1.69%  '
